$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.193.06'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").Value = '1.681.04'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = "'215.91"
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("D6").Value = "'0.5272"
$ws.Range("E6").Value = '  -1.23%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = "'0.2687"
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Value = "'0.06364"
$ws.Range("E9").Value = '  -1.67%  '
$ws.Range("D10").Value = "'21.45"
$ws.Range("E10").Value = '  -2.14%  '
$ws.Range("D11").Value = "'0.07627"
$ws.Range("E11").Value = '  +1.12%  '
$ws.Range("D12").Value = '1.689.84'
$ws.Range("E12").Value = '  +0.25%  '
$ws.Range("D13").Value = "'4.523"
$ws.Range("E13").Value = '  -0.03%  '
$ws.Range("D14").Value = "'0.5748"
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").Value = "'0.000008238"
$ws.Range("E15").Value = '  -2.41%  '
$ws.Range("D16").Value = "'66.41"
$ws.Range("E16").Value = '  +2.54%  '
$ws.Range("D17").Value = '26.232.61'
$ws.Range("E17").Value = '  -0.36%  '
$ws.Range("D18").Value = "'1.006"
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").Value = "'4.867"
$ws.Range("E19").Value = '  -0.79%  '
$ws.Range("D20").Value = "'10.73"
$ws.Range("E20").Value = '  -1.24%  '
$ws.Range("D21").Value = "'189.72"
$ws.Range("E21").Value = '  -0.71%  '
$ws.Range("D22").Value = "'6.231"
$ws.Range("E22").Value = '  +0.38%  '
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").Value = "'149.07"
$ws.Range("E24").Value = '  +2.01%  '
$ws.Range("D25").Value = "'0.1260"
$ws.Range("E25").Value = '  -1.51%  '
$ws.Range("D26").Value = "'7.714"
$ws.Range("E26").Value = '  -1.52%  '
$ws.Range("D27").Value = "'15.83"
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("D28").Value = "'0.06401"
$ws.Range("E28").Value = '  -1.45%  '
$ws.Range("D29").Value = "'1.376"
$ws.Range("E29").Value = '  -0.75%  '
$ws.Range("D30").Value = "'1.314"
$ws.Range("E30").Value = '  -0.68%  '
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("D32").Value = "'3.564"
$ws.Range("E32").Value = '  -0.37%  '
$ws.Range("D33").Value = "'1.681"
$ws.Range("E33").Value = '  +0.84%  '
$ws.Range("D34").Value = "'1.018"
$ws.Range("E34").Value = '  -1.43%  '
$ws.Range("D35").Value = "'0.6114"
$ws.Range("E35").Value = '  -0.81%  '
$ws.Range("D36").Value = "'2.421"
$ws.Range("E36").Value = '  +0.83%  '
$ws.Range("D37").Value = "'2.746"
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = "'0.01638"
$ws.Range("E38").Value = '  +1.07%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = "'6.171"
$ws.Range("E39").Value = '  -1.13%  '
$ws.Range("D40").Value = '1.094.97'
$ws.Range("E40").Value = '  -1.30%  '
$ws.Range("D41").Value = "'0.8821"
$ws.Range("E41").Value = '  +1.39%  '
$ws.Range("E42").Value = '  -0.48%  '
$ws.Range("D43").Value = "'100.42"
$ws.Range("E43").Value = '  -0.12%  '
$ws.Range("D44").Value = '1.835.16'
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = "'57.49"
$ws.Range("E45").Value = '  +0.71%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = "'0.00000000108"
$ws.Range("E46").Value = '  +1.36%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'8.113"
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").Value = "'1.001"
$ws.Range("E48").Value = '  +0.09%  '
$ws.Range("D49").Value = "'0.05266"
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").Value = "'0.4279"
$ws.Range("E50").Value = '  -0.27%  '
$ws.Range("D51").Value = "'6.016"
$ws.Range("E51").Value = '  -0.98%  '
